# Insert a new first column "No" into the Parts sheet, shifting all
# existing headers one column to the right (A->B, B->C, ... O->P).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parts")

# Insert a new column before column A; this shifts the existing data
# (and column formatting) one column to the right automatically.
$ws.Range("A1").EntireColumn.Insert()

# New column A holds the "No" header text.
$ws.Range("A1").Value = "No"

# Narrow width for the new "No" column, centered content for future rows
# (mirrors the existing narrow numeric-style columns in this sheet).
$ws.Columns.Item(1).ColumnWidth = 2.67
$ws.Columns.Item(1).HorizontalAlignment = -4108   # xlCenter
$ws.Columns.Item(1).VerticalAlignment = -4108     # xlCenter
$ws.Columns.Item(1).WrapText = $true

# Re-apply the bold header look to A1 itself (column-level alignment above
# would otherwise also affect the header cell's alignment).
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Update the selection to match the target state (whole column A selected).
$ws.Range("A1:A1048576").Select()
